# Finalized and added the last stuff: mark all CRUD tests as passed
# on the "Test Results" sheet (all FALSE -> TRUE for B2:E24).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

$ws.Range("B2:E24").Value = $true
